$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly logged timesheet hours (rows 10-13) ---
$ws.Range("C10").Value = 11
$ws.Range("E10").Value = 11

$ws.Range("C11").Value = 8
$ws.Range("E11").Value = 7

$ws.Range("C12").Value = 5
$ws.Range("E12").Value = 9

$ws.Range("B13").Value = 8
$ws.Range("C13").Value = 4
$ws.Range("E13").Value = 9

# --- Leave the final selection where the author left off editing ---
[void]$ws.Range("H2").Select()
